# Refresh the crypto-price sheet: overwrite the Price (D) and
# Volume(1h) (E) inline-string cells with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.477.76"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.625.74"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'594.19"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").Value = "'167.68"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("D9").Value = "2.625.64"
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("E10").Value = "  -2.62%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D15").Value = "3.117.59"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").Value = "67.565.63"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "2.628.61"
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("D20").Value = "'8.02"
$ws.Range("E20").Value = "  +2.54%  "
$ws.Range("D21").Value = "'358.02"
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("E24").Value = "  -3.57%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'69.94"
$ws.Range("E27").Value = "  -1.14%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").Value = "'546.76"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("D34").Value = "'1.90"
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("E35").Value = "  +5.25%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E37").Value = "  -2.43%  "
$ws.Range("D38").Value = "'157.87"
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("E45").Value = "  -3.00%  "
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("D47").Value = "'152.84"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("E49").Value = "  -1.43%  "
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("E51").Value = "  -0.66%  "
